$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preserving its General/default
# number format. Values that Excel would otherwise auto-parse as a plain
# number (single decimal point) need the cell pre-formatted as Text so the
# literal (e.g. leading zeros, trailing ".00") is preserved verbatim, exactly
# as the source diff requires (all Price/Volume cells are inline strings).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

$ws.Range("D2").Value = '71.807.07'
$ws.Range("E2").Value = '  +4.32%  '
$ws.Range("D3").Value = '2.632.90'
$ws.Range("E3").Value = '  +4.58%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue $ws.Range("D5") '606.77'
$ws.Range("E5").Value = '  +2.22%  '
Set-TextValue $ws.Range("D6") '179.47'
$ws.Range("E6").Value = '  +2.81%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +1.71%  '
$ws.Range("D9").Value = '2.632.67'
$ws.Range("E9").Value = '  +4.64%  '
Set-TextValue $ws.Range("D10") '0.169'
$ws.Range("E10").Value = '  +13.30%  '
Set-TextValue $ws.Range("D11") '0.165'
$ws.Range("E11").Value = '  +0.41%  '
$ws.Range("E12").Value = '  +3.14%  '
Set-TextValue $ws.Range("D13") '5.05'
$ws.Range("E13").Value = '  +1.46%  '
$ws.Range("D14").Value = '3.134.53'
Set-TextValue $ws.Range("D15") '0.0000187'
$ws.Range("E15").Value = '  +8.79%  '
Set-TextValue $ws.Range("D16") '26.57'
$ws.Range("E16").Value = '  +2.73%  '
$ws.Range("D17").Value = '71.656.77'
$ws.Range("E17").Value = '  +4.26%  '
$ws.Range("D18").Value = '2.629.34'
$ws.Range("E18").Value = '  +4.18%  '
Set-TextValue $ws.Range("D19") '383.19'
$ws.Range("E19").Value = '  +5.41%  '
Set-TextValue $ws.Range("D20") '7.97'
$ws.Range("E20").Value = '  +6.05%  '
Set-TextValue $ws.Range("D21") '11.50'
$ws.Range("E21").Value = '  +4.97%  '
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("E23").Value = '  +18.08%  '
Set-TextValue $ws.Range("D24") '72.75'
$ws.Range("E24").Value = '  +3.21%  '
Set-TextValue $ws.Range("D25") '4.47'
$ws.Range("E25").Value = '  +6.86%  '
$ws.Range("E26").Value = '  +0.07%  '
Set-TextValue $ws.Range("D27") '9.95'
$ws.Range("E27").Value = '  +10.65%  '
$ws.Range("D28").Value = '2.766.99'
$ws.Range("E28").Value = '  +4.60%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").Value = '0.0₃0965'
$ws.Range("E30").Value = '  +9.89%  '
Set-TextValue $ws.Range("D31") '545.72'
$ws.Range("E31").Value = '  +6.97%  '
$ws.Range("E32").Value = '  +3.74%  '
$ws.Range("E33").Value = '  +8.21%  '
$ws.Range("E34").Value = '  +3.20%  '
Set-TextValue $ws.Range("D35") '1.00'
$ws.Range("E35").Value = '  +0.02%  '
Set-TextValue $ws.Range("D36") '166.02'
$ws.Range("E36").Value = '  +2.25%  '
Set-TextValue $ws.Range("D37") '19.24'
$ws.Range("E37").Value = '  +3.21%  '
$ws.Range("E38").Value = '  -2.82%  '
Set-TextValue $ws.Range("D39") '19.15'
$ws.Range("E39").Value = '  +2.63%  '
$ws.Range("E40").Value = '  +6.48%  '
Set-TextValue $ws.Range("D41") '1.87'
$ws.Range("E41").Value = '  +8.04%  '
Set-TextValue $ws.Range("D42") '2.64'
$ws.Range("E42").Value = '  +11.81%  '
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("E44").Value = '  +5.33%  '
Set-TextValue $ws.Range("D45") '0.333'
$ws.Range("E45").Value = '  +2.66%  '
Set-TextValue $ws.Range("D46") '39.24'
$ws.Range("E46").Value = '  +0.60%  '
Set-TextValue $ws.Range("D47") '150.93'
$ws.Range("E47").Value = '  +0.17%  '
Set-TextValue $ws.Range("D48") '3.65'
$ws.Range("E48").Value = '  +2.17%  '
Set-TextValue $ws.Range("D49") '0.536'
$ws.Range("E49").Value = '  +4.23%  '
$ws.Range("E50").Value = '  +7.72%  '
$ws.Range("D51").Value = '0.0₆0265'
$ws.Range("E51").Value = '  +5.69%  '
